$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update priority for "Recibir y responder mensaje del cliente" story (row 6)
$ws.Range("F6").Value = "Baja"

# Update comentarios / elementos de la interfaz for row 6 (whatsapp contact feature)
$ws.Range("H6").Value = "Activity con enlace directo al whatsapp del proveedor"
$ws.Range("I6").Value = "3 Text View con dirección nombre y telefono del proveedor , logo de whatsapp que redirecciona a la app para chatear"

# Row 6 grew taller to fit the new wrapped text
$ws.Rows.Item(6).RowHeight = 82.8

# Mark layout as done for row 8 (Recibir solicitud del cliente)
$ws.Range("H8").Value = "Layout OK"

# Scroll / selection state as left by the author
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("I8").Select()
